# Update countries & provincias Spain
# - Refreshes the COVID-19 snapshot numbers for a handful of countries
#   (the refreshed case counts also change the descending sort order of a
#   few rows that are tied/near the updated ones, but their own data is
#   unchanged).
# - Bumps the "Datos actualizados" footer timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- España (row 5) ---
$ws.Range("B5").Value = 262783
$ws.Range("C5").Value = 2666
$ws.Range("D5").Value = 173157
$ws.Range("E5").Value = 63148
$ws.Range("F5").Value = 1741
$ws.Range("G5").Value = 179
$ws.Range("H5").Value = 26478

# --- Suiza (row 22) ---
$ws.Range("B22").Value = 30251
$ws.Range("C22").Value = 44
$ws.Range("E22").Value = 2328

# --- Rumania (row 37) ---
$ws.Range("B37").Value = 15131
$ws.Range("C37").Value = 320
$ws.Range("D37").Value = 6912
$ws.Range("E37").Value = 7293
$ws.Range("F37").Value = 245

# --- Australia (row 53) ---
$ws.Range("B53").Value = 6929
$ws.Range("C53").Value = 15
$ws.Range("D53").Value = 6135
$ws.Range("E53").Value = 697
$ws.Range("F53").Value = 19

# --- Marruecos (row 56) ---
$ws.Range("B56").Value = 5873
$ws.Range("C56").Value = 162
$ws.Range("D56").Value = 2389
$ws.Range("E56").Value = 3298

# --- Moldavia (row 60) ---
$ws.Range("D60").Value = 1925
$ws.Range("E60").Value = 2644
$ws.Range("G60").Value = 9
$ws.Range("H60").Value = 159

# --- Eslovenia (row 89) ---
$ws.Range("B89").Value = 1454
$ws.Range("C89").Value = 4
$ws.Range("D89").Value = 255
$ws.Range("F89").Value = 10
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 101

# --- Albania gains cases and overtakes Mayotte -> rows 101/102 swap ---
$ws.Range("A101").Value = "Albania"
$ws.Range("B101").Value = 856
$ws.Range("C101").Value = 6
$ws.Range("D101").Value = 627
$ws.Range("E101").Value = 198
$ws.Range("F101").Value = 7
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 31

$ws.Range("A102").Value = "Mayotte"
$ws.Range("B102").Value = 854
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 352
$ws.Range("E102").Value = 492
$ws.Range("F102").Value = 7
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 10

# --- Etiopia gains cases and overtakes Santo Tome y Principe / Liberia -> rows 140/141/142 shift ---
$ws.Range("A140").Value = "Etiopia"
$ws.Range("B140").Value = 210
$ws.Range("C140").Value = 16
$ws.Range("D140").Value = 97
$ws.Range("E140").Value = 108
$ws.Range("F140").Value = 1
$ws.Range("G140").Value = 1
$ws.Range("H140").Value = 5

$ws.Range("A141").Value = "Santo Tome y Principe"
$ws.Range("B141").Value = 208
$ws.Range("C141").Value = 0
$ws.Range("D141").Value = 4
$ws.Range("E141").Value = 199
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 5

$ws.Range("A142").Value = "Liberia"
$ws.Range("B142").Value = 199
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 79
$ws.Range("E142").Value = 100
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 20

# --- Laos (row 189) ---
$ws.Range("D189").Value = 13
$ws.Range("E189").Value = 6

# --- Butan / Islas Virgenes Britanicas swap (rows 212/213) ---
$ws.Range("A212").Value = "Butan"
$ws.Range("B212").Value = 7
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 5
$ws.Range("E212").Value = 2
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 0

$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("B213").Value = 7
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 4
$ws.Range("E213").Value = 2
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

# --- Footer timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 13:04"
